$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tabla de datos dispersos")

# The "P(90) - P(10)" column (L) is removed entirely; this shifts every
# column to its right (Moda, Opinion del 1 al 6, fi, hi, Fi, Hi -- old M:R)
# one position to the left (new L:Q), updating mergeCells/dimension along
# the way, exactly like a user right-clicking column L and choosing Delete.
$ws.Columns("L").Delete()

# A few header labels were shortened/reworded.
$ws.Range("C1").Value = "Nº Valoraciones"
$ws.Range("F1").Value = "Cu"
$ws.Range("K1").Value = "Q3-Q1"

# Kurtosis (F) and skewness (G) were recomputed using the moment-based
# formulas for both book groups (row 2 = Libro 68, row 8 = Libro 69).
$ws.Range("F2").Value = 3.341999482007184
$ws.Range("G2").Value = -1.943931833238586
$ws.Range("F8").Value = 12.74134387819694
$ws.Range("G8").Value = -3.36340775681591
